$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = 1.42
$ws.Range("E2").Value = 1.38

# Row 3
$ws.Range("B3").Value = 1.39
$ws.Range("D3").Value = 1.32
$ws.Range("F3").Value = 1

# Row 4
$ws.Range("C4").Value = 1.54
$ws.Range("E4").Value = 1.15
$ws.Range("F4").Value = 1.03

# Row 5
$ws.Range("B5").Value = 1.38
$ws.Range("C5").Value = 1.3
$ws.Range("D5").Value = 1.39
$ws.Range("E5").Value = 1.22
$ws.Range("F5").Value = 1.08
$ws.Range("G5").Value = 0.58

# Row 6
$ws.Range("C6").Value = 1.64
$ws.Range("D6").Value = 1.62
$ws.Range("E6").Value = 1.28

# Row 7
$ws.Range("E7").Value = 2.08
$ws.Range("F7").Value = 1.47

$wb.Save()
